$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "229.83"
Set-TextValue $ws "G2" "5"
Set-TextValue $ws "D3" "22.32"
Set-TextValue $ws "G3" "5"
Set-TextValue $ws "D4" "5.252"
Set-TextValue $ws "G4" "5"
Set-TextValue $ws "D5" "0.05552"
Set-TextValue $ws "G5" "5"
Set-TextValue $ws "D6" "3.380"
Set-TextValue $ws "G6" "5"
Set-TextValue $ws "D7" "6.470"
Set-TextValue $ws "G7" "5"
Set-TextValue $ws "G8" "5"
Set-TextValue $ws "D9" "0.7784"
Set-TextValue $ws "G9" "5"
Set-TextValue $ws "D10" "0.1380"
Set-TextValue $ws "G10" "5"
Set-TextValue $ws "D11" "0.07326"
Set-TextValue $ws "G11" "5"
Set-TextValue $ws "D12" "0.03134"
Set-TextValue $ws "G12" "5"
Set-TextValue $ws "D13" "0.02946"
Set-TextValue $ws "G13" "5"
Set-TextValue $ws "D14" "0.09267"
Set-TextValue $ws "G14" "5"
Set-TextValue $ws "D15" "0.001667"
Set-TextValue $ws "G15" "5"
Set-TextValue $ws "G16" "5"
Set-TextValue $ws "D17" "0.04786"
Set-TextValue $ws "G17" "5"
Set-TextValue $ws "D18" "0.0005889"
Set-TextValue $ws "G18" "5"
Set-TextValue $ws "D19" "0.006197"
Set-TextValue $ws "G19" "5"
Set-TextValue $ws "D20" "0.005235"
Set-TextValue $ws "G20" "5"
Set-TextValue $ws "D21" "0.001064"
Set-TextValue $ws "G21" "5"
Set-TextValue $ws "G22" "5"
Set-TextValue $ws "D23" "3.913"
Set-TextValue $ws "G23" "5"
Set-TextValue $ws "D24" "2.147"
Set-TextValue $ws "G24" "5"
Set-TextValue $ws "G25" "5"
Set-TextValue $ws "G26" "5"
$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"
Set-TextValue $ws "G27" "5"
Set-TextValue $ws "G28" "5"
Set-TextValue $ws "G29" "5"
Set-TextValue $ws "G30" "5"
Set-TextValue $ws "G31" "5"
Set-TextValue $ws "G32" "5"
Set-TextValue $ws "G33" "5"
Set-TextValue $ws "G34" "5"
Set-TextValue $ws "G35" "5"
Set-TextValue $ws "G36" "5"
Set-TextValue $ws "G37" "5"
Set-TextValue $ws "G38" "5"
Set-TextValue $ws "G39" "5"
Set-TextValue $ws "D40" "0.03995"
Set-TextValue $ws "G40" "5"
Set-TextValue $ws "D41" "0.007138"
$ws.Range("E41").Value = "40KickTokenKICK"
Set-TextValue $ws "G41" "5"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D42" "0.1040"
$ws.Range("E42").Value = "41BKEXTokenBKK"
Set-TextValue $ws "G42" "5"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D43" "0.002980"
$ws.Range("E43").Value = "42CEJICEJI"
Set-TextValue $ws "G43" "5"
Set-TextValue $ws "D44" "0.009991"
Set-TextValue $ws "G44" "5"
Set-TextValue $ws "D45" "0.00005438"
Set-TextValue $ws "G45" "5"
Set-TextValue $ws "G46" "5"
Set-TextValue $ws "D47" "0.7851"
Set-TextValue $ws "G47" "5"
Set-TextValue $ws "D48" "0.04159"
Set-TextValue $ws "G48" "5"
Set-TextValue $ws "G49" "5"
Set-TextValue $ws "G50" "5"
Set-TextValue $ws "G51" "5"
